$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.100.19'
$ws.Range("E2").Value = '  -1.00%  '
$ws.Range("D3").Value = '3.009.87'
$ws.Range("E3").Value = '  -0.33%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '593.53'
$ws.Range("E5").Value = '  +1.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '146.98'
$ws.Range("E6").Value = '  -1.14%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '3.011.63'
$ws.Range("E8").Value = '  -0.20%  '
$ws.Range("E9").Value = '  -2.04%  '
$ws.Range("E10").Value = '  +8.08%  '
$ws.Range("E11").Value = '  -0.81%  '
$ws.Range("E12").Value = '  -1.18%  '
$ws.Range("E13").Value = '  -0.04%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.35'
$ws.Range("E14").Value = '  -1.63%  '
$ws.Range("E15").Value = '  +2.51%  '
$ws.Range("D16").Value = '3.515.63'
$ws.Range("E16").Value = '  -0.16%  '
$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.98'
$ws.Range("E17").Value = '  -2.14%  '
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '62.120.56'
$ws.Range("E18").Value = '  -0.89%  '
$ws.Range("D19").Value = '3.009.49'
$ws.Range("E19").Value = '  -0.29%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '445.86'
$ws.Range("E20").Value = '  -3.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.16'
$ws.Range("E21").Value = '  +0.71%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.688'
$ws.Range("E22").Value = '  -0.82%  '
$ws.Range("E23").Value = '  -0.92%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.18'
$ws.Range("E24").Value = '  +0.40%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.93'
$ws.Range("E25").Value = '  +8.74%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.25'
$ws.Range("E26").Value = '  +0.69%  '
$ws.Range("E27").Value = '  -2.47%  '
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("E29").Value = '  +1.53%  '
$ws.Range("E30").Value = '  +0.15%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.19'
$ws.Range("E31").Value = '  +1.63%  '
$ws.Range("E32").Value = '  -0.67%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.41'
$ws.Range("E33").Value = '  -2.62%  '
$ws.Range("E34").Value = '  +0.30%  '
$ws.Range("E35").Value = '  +3.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.02'
$ws.Range("E36").Value = '  -0.93%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.83'
$ws.Range("E37").Value = '  +0.35%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '50.24'
$ws.Range("E38").Value = '  -0.44%  '
$ws.Range("E39").Value = '  -4.60%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '9.03'
$ws.Range("E40").Value = '  -1.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.95'
$ws.Range("E41").Value = '  +0.25%  '
$ws.Range("E42").Value = '  +0.89%  '
$ws.Range("E43").Value = '  +11.18%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.284'
$ws.Range("E44").Value = '  +4.73%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '394.90'
$ws.Range("E45").Value = '  -0.29%  '
$ws.Range("E46").Value = '  -2.73%  '
$ws.Range("D47").Value = '2.725.29'
$ws.Range("E47").Value = '  -0.77%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '134.12'
$ws.Range("E48").Value = '  +3.69%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.17'
$ws.Range("E50").Value = '  -1.75%  '
$ws.Range("E51").Value = '  -1.88%  '
